# Fruta / hortaliza, semanal
# Update the weekly data block (rows 2-19) on the active sheet:
# columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ D = 44508; J = 400;  K = 900;  L = 1000; M = 950;  P = 950  }
    3  = @{ D = 44525; J = 360;  K = 800;  L = 900;  M = 850;  P = 850  }
    4  = @{ D = 44518; J = 400;  K = 800;  L = 900;  M = 850;  P = 850  }
    5  = @{ D = 44532; J = 240;  K = 800;  L = 900;  M = 850;  P = 850  }
    6  = @{ D = 44503; J = 400;  K = 900;  L = 1000; M = 950;  P = 950  }
    7  = @{ D = 44505; J = 440;  K = 900;  L = 1000; M = 950;  P = 950  }
    8  = @{ D = 44523; J = 400;  K = 800;  L = 900;  M = 850;  P = 850  }
    9  = @{ D = 44510; J = 600;  K = 900;  L = 1000; M = 950;  P = 950  }
    10 = @{ D = 44511; J = 500;  K = 900;  L = 1000; M = 950;  P = 950  }
    11 = @{ D = 44545; J = 4000; K = 800;  L = 900;  M = 850;  P = 850  }
    12 = @{ D = 44530; J = 300;  K = 800;  L = 900;  M = 850;  P = 850  }
    13 = @{ D = 44553; J = 8000; K = 800;  L = 900;  M = 850;  P = 850  }
    14 = @{ D = 44516; J = 400;  K = 900;  L = 1000; M = 950;  P = 950  }
    15 = @{ D = 44476; J = 300;  K = 1100; L = 1200; M = 1150; P = 1150 }
    16 = @{ D = 44524; J = 400;  K = 800;  L = 900;  M = 850;  P = 850  }
    17 = @{ D = 44517; J = 500;  K = 800;  L = 900;  M = 850;  P = 850  }
    18 = @{ D = 44537; J = 400;  K = 800;  L = 900;  M = 850;  P = 850  }
    19 = @{ D = 44504; J = 500;  K = 900;  L = 1000; M = 950;  P = 950  }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
